$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data values of row 2 and row 3 in columns
# A, B, E, F, G, H, Q, R (all other columns are identical between
# the two rows so no visible change occurs there).
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $val2 = $ws.Range($addr2).Value()
    $val3 = $ws.Range($addr3).Value()
    $ws.Range($addr2).Value = $val3
    $ws.Range($addr3).Value = $val2
}
